$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "ESPETOS"
$ws.Range("C3").Value = "ESPETOS"
$ws.Range("A4").Value = "Linguiça Caseira"
$ws.Range("C4").Value = "ESPETOS"
$ws.Range("C5").Value = "ESPETOS"
$ws.Range("C6").Value = "ESPETOS"
$ws.Range("A7").Value = "Asa de Frango"
$ws.Range("C7").Value = "ESPETOS"
$ws.Range("C8").Value = "ESPETOS"
$ws.Range("A9").Value = "Linguiça Mista"
$ws.Range("C9").Value = "ESPETOS"
$ws.Range("C10").Value = "ESPETOS"
$ws.Range("A11").Value = "Coração de Galinha"
$ws.Range("C11").Value = "ESPETOS"
$ws.Range("C12").Value = "ESPETOS"
$ws.Range("C13").Value = "ESPETOS"
$ws.Range("A14").Value = "Cupim c/ Queijo"
$ws.Range("C14").Value = "ESPETOS"
$ws.Range("G14").Value = 9
$ws.Range("A15").Value = "Baião Enxuto P"
$ws.Range("C15").Value = "ACOMPANHAMENTOS"
$ws.Range("G15").Value = 12
$ws.Range("A16").Value = "Baião Enxuto G"
$ws.Range("C16").Value = "ACOMPANHAMENTOS"
$ws.Range("G16").Value = 15
$ws.Range("A17").Value = "Baião Cremoso P"
$ws.Range("C17").Value = "ACOMPANHAMENTOS"
$ws.Range("G17").Value = 15
$ws.Range("A18").Value = "Baião Cremoso G"
$ws.Range("C18").Value = "ACOMPANHAMENTOS"
$ws.Range("G18").Value = 18
$ws.Range("A19").Value = "Arroz P"
$ws.Range("C19").Value = "ACOMPANHAMENTOS"
$ws.Range("G19").Value = 12
$ws.Range("A20").Value = "Arroz G"
$ws.Range("C20").Value = "ACOMPANHAMENTOS"
$ws.Range("G20").Value = 14
$ws.Range("A21").Value = "Macaxeira Frita"
$ws.Range("C21").Value = "ACOMPANHAMENTOS"
$ws.Range("G21").Value = 15
$ws.Range("A22").Value = "Creme de Macaxeira"
$ws.Range("C22").Value = "ACOMPANHAMENTOS"
$ws.Range("G22").Value = 15
$ws.Range("A23").Value = "Batata Frita"
$ws.Range("C23").Value = "ACOMPANHAMENTOS"
$ws.Range("G23").Value = 17
$ws.Range("A24").Value = "Piabinha"
$ws.Range("C24").Value = "ACOMPANHAMENTOS"
$ws.Range("G24").Value = 17
$ws.Range("A25").Value = "Torresmo"
$ws.Range("C25").Value = "ACOMPANHAMENTOS"
$ws.Range("A26").Value = "Caldinho de Feijão"
$ws.Range("C26").Value = "ACOMPANHAMENTOS"
$ws.Range("G26").Value = 12
$ws.Range("A27").Value = "Camarão c/ Alho e Óleo"
$ws.Range("C27").Value = "ACOMPANHAMENTOS"
$ws.Range("G27").Value = 30
$ws.Range("A28").Value = "Queijo na chapa"
$ws.Range("C28").Value = "ACOMPANHAMENTOS"
$ws.Range("G28").Value = 15
$ws.Range("A29").Value = "Carne de sol c/ Macaxeira"
$ws.Range("C29").Value = "ACOMPANHAMENTOS"
$ws.Range("G29").Value = 30
$ws.Range("C30").Value = "ACOMPANHAMENTOS"
$ws.Range("C31").Value = "ACOMPANHAMENTOS"
$ws.Range("C32").Value = "JANTINHAS"
$ws.Range("C33").Value = "JANTINHAS"
$ws.Range("C34").Value = "JANTINHAS"
$ws.Range("A35").Value = "Porção Mungunzá mesa"
$ws.Range("C35").Value = "MUNGUNZÁ"
$ws.Range("G35").Value = 10
$ws.Range("A36").Value = "Porção Mungunzá viagem"
$ws.Range("C36").Value = "MUNGUNZÁ"
$ws.Range("G36").Value = 12
$ws.Range("A37").Value = "Refri Lata"
$ws.Range("C37").Value = "BEBIDAS"
$ws.Range("G37").Value = 5
$ws.Range("A38").Value = "Refri 600ml"
$ws.Range("C38").Value = "BEBIDAS"
$ws.Range("G38").Value = 9
$ws.Range("A39").Value = "Refri 1L"
$ws.Range("C39").Value = "BEBIDAS"
$ws.Range("G39").Value = 10
$ws.Range("A40").Value = "Refri 1,5L"
$ws.Range("C40").Value = "BEBIDAS"
$ws.Range("G40").Value = 11
$ws.Range("A41").Value = "Água"
$ws.Range("C41").Value = "BEBIDAS"
$ws.Range("G41").Value = 3
$ws.Range("A42").Value = "Água c/ gás"
$ws.Range("C42").Value = "BEBIDAS"
$ws.Range("G42").Value = 3.5
$ws.Range("A43").Value = "Cerveja Buchuda"
$ws.Range("C43").Value = "BEBIDAS"
$ws.Range("G43").Value = 4.5
$ws.Range("A44").Value = "Long Neck"
$ws.Range("C44").Value = "BEBIDAS"
$ws.Range("G44").Value = 8
$ws.Range("A45").Value = "Dose"
$ws.Range("C45").Value = "BEBIDAS"
$ws.Range("G45").Value = 5
$ws.Range("A46").Value = "Copo"
$ws.Range("C46").Value = "BEBIDAS"
$ws.Range("G46").Value = 10
$ws.Range("A47").Value = "Suco S/L"
$ws.Range("C47").Value = "BEBIDAS"
$ws.Range("G47").Value = 6
$ws.Range("A48").Value = "Suco C/L"
$ws.Range("C48").Value = "BEBIDAS"
$ws.Range("G48").Value = 7
$ws.Range("A49").Value = "Cerveja 600ml"
$ws.Range("C49").Value = "BEBIDAS"
$ws.Range("G49").Value = 14
